$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 80, shifting existing rows 80+ down by one.
$ws.Rows("80:80").Insert()

# Populate the newly inserted row 80 with the new price record.
$ws.Range("A80").Value = 3
$ws.Range("B80").Value = "Femacal de La Calera"
$ws.Range("C80").Value = "Coquimbo"
$ws.Range("D80").Value = 44601
$ws.Range("E80").Value = 5
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100101
$ws.Range("H80").Value = "Berries"
$ws.Range("I80").Value = 100101001
$ws.Range("J80").Value = "Arándano (blue)"
$ws.Range("K80").Value = "Sin especificar"
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 45
$ws.Range("N80").Value = 4000
$ws.Range("O80").Value = 4500
$ws.Range("P80").Value = 4222
$ws.Range("Q80").Value = "`$/bandeja 2 kilos"
$ws.Range("R80").Value = "Provincia de Linares"
$ws.Range("S80").Value = 2111
$ws.Range("T80").Value = 2
